# Scheduled market-data refresh: update the computed price/profit columns
# (H:N -> currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) on each class sheet with the latest snapshot values.
# Some rows gain/lose their LeveProfitNQ (M) or LeveProfitHQ (N) cell
# entirely depending on whether that market is currently populated, so
# ClearContents() is used (rather than writing 0/blank) to drop a column
# that should no longer be present for that row.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1012.9231
$ws.Range("I107").Value = 605.7778
$ws.Range("K107").Value = 605.7778
$ws.Range("M107").Value = 1314.2222
$ws.Range("H132").Value = 6040.081
$ws.Range("I132").Value = 6376.5757
$ws.Range("K132").Value = 19129.7271
$ws.Range("M132").Value = -16599.7271

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1336
$ws.Range("I2").Value = 1068.6364
$ws.Range("K2").Value = 1068.6364
$ws.Range("M2").Value = -955.6364000000001
$ws.Range("H12").Value = 900
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 900
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = 900
$ws.Range("N12").Value = -1246
$ws.Range("H32").Value = 18958.887
$ws.Range("I32").Value = 18733.846
$ws.Range("K32").Value = 18733.846
$ws.Range("M32").Value = -18446.846
$ws.Range("H110").Value = 3600.4443
$ws.Range("I110").Value = 3401.6667
$ws.Range("K110").Value = 3401.6667
$ws.Range("M110").Value = -1356.6667
$ws.Range("H116").Value = 1336
$ws.Range("I116").Value = 1068.6364
$ws.Range("K116").Value = 1068.6364
$ws.Range("M116").Value = 1225.3636
$ws.Range("H132").Value = 59969.89
$ws.Range("I132").Value = 62791.65
$ws.Range("K132").Value = 188374.95
$ws.Range("M132").Value = -185844.95
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("N139").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("L139").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1336
$ws.Range("I3").Value = 1068.6364
$ws.Range("K3").Value = 1068.6364
$ws.Range("M3").Value = -954.6364000000001
$ws.Range("H86").Value = 3569.32
$ws.Range("I86").Value = 1868.091
$ws.Range("K86").Value = 1868.091
$ws.Range("M86").Value = -745.0909999999999
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("H89").Value = 3569.32
$ws.Range("I89").Value = 1868.091
$ws.Range("K89").Value = 9340.455
$ws.Range("M89").Value = -3724.455
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").Value = 0
$ws.Range("H99").Value = 174822.33
$ws.Range("I99").Value = 339661.34
$ws.Range("K99").Value = 339661.34
$ws.Range("M99").Value = -338163.34
$ws.Range("H105").Value = 3343.8262
$ws.Range("I105").Value = 3381.0454
$ws.Range("K105").Value = 3381.0454
$ws.Range("M105").Value = -1634.0454
$ws.Range("H134").Value = 2921.8333
$ws.Range("I134").Value = 2104.6667
$ws.Range("K134").Value = 6314.000100000001
$ws.Range("M134").Value = -3779.000100000001
$ws.Range("L87").ClearContents()
$ws.Range("L90").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 2676.25
$ws.Range("I35").Value = 2676.25
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 2676.25
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = -2382.25
$ws.Range("M35").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1100
$ws.Range("I59").Value = 1100
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 3300
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = -2760
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("M101").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 5668.6665
$ws.Range("I36").Value = 7008
$ws.Range("J36").Value = 2990
$ws.Range("K36").Value = 7008
$ws.Range("L36").Value = 2990
$ws.Range("M36").Value = -6523
$ws.Range("N36").Value = -3960
$ws.Range("H43").Value = 7728.3335
$ws.Range("I43").Value = 4093
$ws.Range("K43").Value = 4093
$ws.Range("M43").Value = -3942
$ws.Range("H70").Value = 7499
$ws.Range("J70").Value = 7499
$ws.Range("L70").Value = 7499
$ws.Range("N70").Value = -8039
$ws.Range("H73").Value = 7499
$ws.Range("J73").Value = 7499
$ws.Range("L73").Value = 7499
$ws.Range("N73").Value = -9371
$ws.Range("H107").Value = 84554.5
$ws.Range("I107").Value = 167394.17
$ws.Range("J107").Value = 1714.8334
$ws.Range("K107").Value = 167394.17
$ws.Range("L107").Value = 1714.8334
$ws.Range("M107").Value = -165474.17
$ws.Range("N107").Value = -5554.8334
$ws.Range("H122").Value = 4660.5454
$ws.Range("I122").Value = 3866
$ws.Range("J122").Value = 4958.5
$ws.Range("K122").Value = 11598
$ws.Range("L122").Value = 14875.5
$ws.Range("M122").Value = -9148
$ws.Range("N122").Value = -19775.5
$ws.Range("H132").Value = 39056.37
$ws.Range("I132").Value = 40466.23
$ws.Range("K132").Value = 121398.69
$ws.Range("M132").Value = -118868.69

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2055.4285
$ws.Range("I16").Value = 1930.2693
$ws.Range("J16").Value = 3682.5
$ws.Range("K16").Value = 1930.2693
$ws.Range("L16").Value = 3682.5
$ws.Range("M16").Value = -1760.2693
$ws.Range("N16").Value = -4022.5
$ws.Range("H36").Value = 69999.5
$ws.Range("J36").Value = 69999.5
$ws.Range("L36").Value = 69999.5
$ws.Range("N36").Value = -71123.5
$ws.Range("H68").Value = 5332.9375
$ws.Range("I68").Value = 3505.25
$ws.Range("K68").Value = 3505.25
$ws.Range("M68").Value = -2756.25
$ws.Range("H71").Value = 5332.9375
$ws.Range("I71").Value = 3505.25
$ws.Range("K71").Value = 17526.25
$ws.Range("M71").Value = -13782.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("H41").Value = 32950
$ws.Range("I41").Value = 32950
$ws.Range("K41").Value = 32950
$ws.Range("M41").Value = -32560
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("N95").Value = 0
$ws.Range("H132").Value = 52626.9
$ws.Range("I132").Value = 58307.39
$ws.Range("K132").Value = 174922.17
$ws.Range("M132").Value = -172392.17
$ws.Range("L16").ClearContents()
$ws.Range("L95").ClearContents()
